$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data (fecha serial 44484) is inserted at the top of the
# data block (row 23), pushing the existing rows down by three and
# leaving the previously-last three rows of data to land at rows 39-41.
$ws.Rows("23:25").Insert()

# Row 23: Espárragos, Sin especificar, Banquete
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44484
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 300000000
$ws.Cells.Item(23, 7).Value = "Espárragos"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Banquete"
$ws.Cells.Item(23, 10).Value = 770
$ws.Cells.Item(23, 11).Value = 1400
$ws.Cells.Item(23, 12).Value = 1500
$ws.Cells.Item(23, 13).Value = 1455
$ws.Cells.Item(23, 14).Value = "$/kilo"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1455
$ws.Cells.Item(23, 17).Value = 1
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Row 24: Espárragos, Sin especificar, Primera
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44484
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = 300000000
$ws.Cells.Item(24, 7).Value = "Espárragos"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 630
$ws.Cells.Item(24, 11).Value = 1200
$ws.Cells.Item(24, 12).Value = 1300
$ws.Cells.Item(24, 13).Value = 1260
$ws.Cells.Item(24, 14).Value = "$/kilo"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 1260
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"

# Row 25: Espárragos, Sin especificar, Segunda
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44484
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 300000000
$ws.Cells.Item(25, 7).Value = "Espárragos"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Segunda"
$ws.Cells.Item(25, 10).Value = 410
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 1100
$ws.Cells.Item(25, 13).Value = 1063
$ws.Cells.Item(25, 14).Value = "$/kilo"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 1063
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"
